$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column D: "canonical SMILES" (canonical, non-isomeric SMILES - i.e.
# the "canonical isomeric SMILES" in column C with stereo-bond markers
# ('/' and '\') stripped out).
$ws.Range("D2").Value = "canonical SMILES"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row() + $usedRange.Rows.Count() - 1
for ($r = 3; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3).Value()
    if ($c -ne $null) {
        $d = $c.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $d
    }
}

# Match the formatting of the corresponding C cell for each row in column D
# (direct ".Style =" assignment isn't reliably applied by this runtime, so
# copy the individual format properties instead).
for ($r = 2; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $dst.Font.Bold = $src.Font.Bold()
    $dst.Font.Size = $src.Font.Size()
    $dst.Interior.Color = $src.Interior.Color()
    $dst.HorizontalAlignment = $src.HorizontalAlignment()
    $dst.VerticalAlignment = $src.VerticalAlignment()
    $dst.WrapText = $src.WrapText()
    $dst.Borders.LineStyle = $src.Borders.LineStyle()
    $dst.Borders.Color = $src.Borders.Color()
    if ($r -ge 3) {
        $dst.ShrinkToFit = $true
    }
}

# New column width for column D (target OOXML width ~= 36.85546875 character
# units; ColumnWidth = 36 is the closest value this runtime's pixel-rounding
# can reach)
$ws.Columns.Item(4).ColumnWidth = 36
